$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 48.68074166666667
$ws.Range("H2").Value = 146.042225
$ws.Range("I2").Value = 0.1601707305796784
$ws.Range("J2").Value = 0.1601707305796784
$ws.Range("M2").Value = 2.526761
$ws.Range("N2").Value = 7.580283
$ws.Range("O2").Value = 0.09514789101715561
$ws.Range("P2").Value = 0.09514789101715561
$ws.Range("Q2").Value = 123.0045994944083
$ws.Range("R2").Value = 1107.041395449675
$ws.Range("S2").Value = 0.01523990721733344
$ws.Range("T2").Value = 0.01523990721733344
$ws.Range("G3").Value = 48.68074166666667
$ws.Range("H3").Value = 146.042225
$ws.Range("I3").Value = 0.1601707305796784
$ws.Range("J3").Value = 0.1601707305796784
$ws.Range("O3").Value = 0.4908680200032562
$ws.Range("P3").Value = 0.4908680200032562
$ws.Range("Q3").Value = 634.580793747989
$ws.Range("R3").Value = 5711.2271437319
$ws.Range("S3").Value = 0.07862268938212176
$ws.Range("T3").Value = 0.07862268938212176
$ws.Range("G4").Value = 48.68074166666667
$ws.Range("H4").Value = 146.042225
$ws.Range("I4").Value = 0.1601707305796784
$ws.Range("J4").Value = 0.1601707305796784
$ws.Range("M4").Value = 2.743855666666667
$ws.Range("N4").Value = 8.231567
$ws.Range("O4").Value = 0.1033228231474227
$ws.Range("P4").Value = 0.1033228231474227
$ws.Range("Q4").Value = 133.5729288796195
$ws.Range("R4").Value = 1202.156359916575
$ws.Range("S4").Value = 0.0165492920690776
$ws.Range("T4").Value = 0.0165492920690776
$ws.Range("G5").Value = 48.68074166666667
$ws.Range("H5").Value = 146.042225
$ws.Range("I5").Value = 0.1601707305796784
$ws.Range("J5").Value = 0.1601707305796784
$ws.Range("M5").Value = 5.865492
$ws.Range("N5").Value = 17.596476
$ws.Range("O5").Value = 0.2208713818117337
$ws.Range("P5").Value = 0.2208713818117337
$ws.Range("Q5").Value = 285.5365007999
$ws.Range("R5").Value = 2569.8285071991
$ws.Range("S5").Value = 0.03537713058892848
$ws.Range("T5").Value = 0.03537713058892849
$ws.Range("G6").Value = 48.68074166666667
$ws.Range("H6").Value = 146.042225
$ws.Range("I6").Value = 0.1601707305796784
$ws.Range("J6").Value = 0.1601707305796784
$ws.Range("M6").Value = 2.384473
$ws.Range("N6").Value = 7.153419
$ws.Range("O6").Value = 0.08978988402043174
$ws.Range("P6").Value = 0.08978988402043174
$ws.Range("Q6").Value = 116.0779141241417
$ws.Range("R6").Value = 1044.701227117275
$ws.Range("S6").Value = 0.01438171132221715
$ws.Range("T6").Value = 0.01438171132221715
$ws.Range("I7").Value = 0.4466378584920503
$ws.Range("J7").Value = 0.4466378584920503
$ws.Range("M7").Value = 2.526761
$ws.Range("N7").Value = 7.580283
$ws.Range("O7").Value = 0.09514789101715561
$ws.Range("P7").Value = 0.09514789101715561
$ws.Range("Q7").Value = 342.999689793668
$ws.Range("R7").Value = 3086.997208143012
$ws.Range("S7").Value = 0.04249665028393738
$ws.Range("T7").Value = 0.04249665028393737
$ws.Range("I8").Value = 0.4466378584920503
$ws.Range("J8").Value = 0.4466378584920503
$ws.Range("O8").Value = 0.4908680200032562
$ws.Range("P8").Value = 0.4908680200032562
$ws.Range("S8").Value = 0.2192402412564873
$ws.Range("T8").Value = 0.2192402412564873
$ws.Range("I9").Value = 0.4466378584920503
$ws.Range("J9").Value = 0.4466378584920503
$ws.Range("M9").Value = 2.743855666666667
$ws.Range("N9").Value = 8.231567
$ws.Range("O9").Value = 0.1033228231474227
$ws.Range("P9").Value = 0.1033228231474227
$ws.Range("Q9").Value = 372.4695934855987
$ws.Range("R9").Value = 3352.226341370388
$ws.Range("S9").Value = 0.04614788446391771
$ws.Range("T9").Value = 0.04614788446391771
$ws.Range("I10").Value = 0.4466378584920503
$ws.Range("J10").Value = 0.4466378584920503
$ws.Range("M10").Value = 5.865492
$ws.Range("N10").Value = 17.596476
$ws.Range("O10").Value = 0.2208713818117337
$ws.Range("P10").Value = 0.2208713818117337
$ws.Range("Q10").Value = 796.221699039696
$ws.Range("R10").Value = 7165.995291357263
$ws.Range("S10").Value = 0.09864952097457273
$ws.Range("T10").Value = 0.09864952097457273
$ws.Range("I11").Value = 0.4466378584920503
$ws.Range("J11").Value = 0.4466378584920503
$ws.Range("M11").Value = 2.384473
$ws.Range("N11").Value = 7.153419
$ws.Range("O11").Value = 0.08978988402043174
$ws.Range("P11").Value = 0.08978988402043174
$ws.Range("Q11").Value = 323.6845508227241
$ws.Range("R11").Value = 2913.160957404516
$ws.Range("S11").Value = 0.0401035615131352
$ws.Range("T11").Value = 0.0401035615131352
$ws.Range("G12").Value = 44.00775933333333
$ws.Range("H12").Value = 132.023278
$ws.Range("I12").Value = 0.1447955540993982
$ws.Range("J12").Value = 0.1447955540993982
$ws.Range("M12").Value = 2.526761
$ws.Range("N12").Value = 7.580283
$ws.Range("O12").Value = 0.09514789101715561
$ws.Range("P12").Value = 0.09514789101715561
$ws.Range("Q12").Value = 111.1970899808527
$ws.Range("R12").Value = 1000.773809827674
$ws.Range("S12").Value = 0.0137769916012182
$ws.Range("T12").Value = 0.0137769916012182
$ws.Range("G13").Value = 44.00775933333333
$ws.Range("H13").Value = 132.023278
$ws.Range("I13").Value = 0.1447955540993982
$ws.Range("J13").Value = 0.1447955540993982
$ws.Range("O13").Value = 0.4908680200032562
$ws.Range("P13").Value = 0.4908680200032562
$ws.Range("Q13").Value = 573.6658459322391
$ws.Range("R13").Value = 5162.992613390153
$ws.Range("S13").Value = 0.07107550694604597
$ws.Range("T13").Value = 0.07107550694604597
$ws.Range("G14").Value = 44.00775933333333
$ws.Range("H14").Value = 132.023278
$ws.Range("I14").Value = 0.1447955540993982
$ws.Range("J14").Value = 0.1447955540993982
$ws.Range("M14").Value = 2.743855666666667
$ws.Range("N14").Value = 8.231567
$ws.Range("O14").Value = 0.1033228231474227
$ws.Range("P14").Value = 0.1033228231474227
$ws.Range("Q14").Value = 120.7509398240695
$ws.Range("R14").Value = 1086.758458416626
$ws.Range("S14").Value = 0.01496068542874519
$ws.Range("T14").Value = 0.01496068542874519
$ws.Range("G15").Value = 44.00775933333333
$ws.Range("H15").Value = 132.023278
$ws.Range("I15").Value = 0.1447955540993982
$ws.Range("J15").Value = 0.1447955540993982
$ws.Range("M15").Value = 5.865492
$ws.Range("N15").Value = 17.596476
$ws.Range("O15").Value = 0.2208713818117337
$ws.Range("P15").Value = 0.2208713818117337
$ws.Range("Q15").Value = 258.127160307592
$ws.Range("R15").Value = 2323.144442768328
$ws.Range("S15").Value = 0.03198119411412972
$ws.Range("T15").Value = 0.03198119411412972
$ws.Range("G16").Value = 44.00775933333333
$ws.Range("H16").Value = 132.023278
$ws.Range("I16").Value = 0.1447955540993982
$ws.Range("J16").Value = 0.1447955540993982
$ws.Range("M16").Value = 2.384473
$ws.Range("N16").Value = 7.153419
$ws.Range("O16").Value = 0.08978988402043174
$ws.Range("P16").Value = 0.08978988402043174
$ws.Range("Q16").Value = 104.9353139208313
$ws.Range("R16").Value = 944.4178252874821
$ws.Range("S16").Value = 0.01300117600925911
$ws.Range("T16").Value = 0.01300117600925911
$ws.Range("G17").Value = 27.64718166666666
$ws.Range("H17").Value = 82.94154499999999
$ws.Range("I17").Value = 0.09096552629253131
$ws.Range("J17").Value = 0.0909655262925313
$ws.Range("M17").Value = 2.526761
$ws.Range("N17").Value = 7.580283
$ws.Range("O17").Value = 0.09514789101715561
$ws.Range("P17").Value = 0.09514789101715561
$ws.Range("Q17").Value = 69.85782039524832
$ws.Range("R17").Value = 628.7203835572349
$ws.Range("S17").Value = 0.008655177981999973
$ws.Range("T17").Value = 0.008655177981999971
$ws.Range("G18").Value = 27.64718166666666
$ws.Range("H18").Value = 82.94154499999999
$ws.Range("I18").Value = 0.09096552629253131
$ws.Range("J18").Value = 0.0909655262925313
$ws.Range("O18").Value = 0.4908680200032562
$ws.Range("P18").Value = 0.4908680200032562
$ws.Range("Q18").Value = 360.3965323096422
$ws.Range("R18").Value = 3243.56879078678
$ws.Range("S18").Value = 0.04465206777976899
$ws.Range("T18").Value = 0.04465206777976898
$ws.Range("G19").Value = 27.64718166666666
$ws.Range("H19").Value = 82.94154499999999
$ws.Range("I19").Value = 0.09096552629253131
$ws.Range("J19").Value = 0.0909655262925313
$ws.Range("M19").Value = 2.743855666666667
$ws.Range("N19").Value = 8.231567
$ws.Range("O19").Value = 0.1033228231474227
$ws.Range("P19").Value = 0.1033228231474227
$ws.Range("Q19").Value = 75.85987608344611
$ws.Range("R19").Value = 682.7388847510149
$ws.Range("S19").Value = 0.00939881498563544
$ws.Range("T19").Value = 0.00939881498563544
$ws.Range("G20").Value = 27.64718166666666
$ws.Range("H20").Value = 82.94154499999999
$ws.Range("I20").Value = 0.09096552629253131
$ws.Range("J20").Value = 0.0909655262925313
$ws.Range("M20").Value = 5.865492
$ws.Range("N20").Value = 17.596476
$ws.Range("O20").Value = 0.2208713818117337
$ws.Range("P20").Value = 0.2208713818117337
$ws.Range("Q20").Value = 162.16432288838
$ws.Range("R20").Value = 1459.47890599542
$ws.Range("S20").Value = 0.02009168148946298
$ws.Range("T20").Value = 0.02009168148946298
$ws.Range("G21").Value = 27.64718166666666
$ws.Range("H21").Value = 82.94154499999999
$ws.Range("I21").Value = 0.09096552629253131
$ws.Range("J21").Value = 0.0909655262925313
$ws.Range("M21").Value = 2.384473
$ws.Range("N21").Value = 7.153419
$ws.Range("O21").Value = 0.08978988402043174
$ws.Range("P21").Value = 0.08978988402043174
$ws.Range("Q21").Value = 65.92395821026167
$ws.Range("R21").Value = 593.315623892355
$ws.Range("S21").Value = 0.008167784055663921
$ws.Range("T21").Value = 0.008167784055663919
$ws.Range("G22").Value = 47.847851
$ws.Range("H22").Value = 143.543553
$ws.Range("I22").Value = 0.1574303305363417
$ws.Range("J22").Value = 0.1574303305363417
$ws.Range("M22").Value = 2.526761
$ws.Range("N22").Value = 7.580283
$ws.Range("O22").Value = 0.09514789101715561
$ws.Range("P22").Value = 0.09514789101715561
$ws.Range("Q22").Value = 120.900083840611
$ws.Range("R22").Value = 1088.100754565499
$ws.Range("S22").Value = 0.01497916393266663
$ws.Range("T22").Value = 0.01497916393266663
$ws.Range("G23").Value = 47.847851
$ws.Range("H23").Value = 143.543553
$ws.Range("I23").Value = 0.1574303305363417
$ws.Range("J23").Value = 0.1574303305363417
$ws.Range("O23").Value = 0.4908680200032562
$ws.Range("P23").Value = 0.4908680200032562
$ws.Range("Q23").Value = 623.7235963786947
$ws.Range("R23").Value = 5613.512367408252
$ws.Range("S23").Value = 0.07727751463883223
$ws.Range("T23").Value = 0.07727751463883223
$ws.Range("G24").Value = 47.847851
$ws.Range("H24").Value = 143.543553
$ws.Range("I24").Value = 0.1574303305363417
$ws.Range("J24").Value = 0.1574303305363417
$ws.Range("M24").Value = 2.743855666666667
$ws.Range("N24").Value = 8.231567
$ws.Range("O24").Value = 0.1033228231474227
$ws.Range("P24").Value = 0.1033228231474227
$ws.Range("Q24").Value = 131.2875971041723
$ws.Range("R24").Value = 1181.588373937551
$ws.Range("S24").Value = 0.01626614620004673
$ws.Range("T24").Value = 0.01626614620004673
$ws.Range("G25").Value = 47.847851
$ws.Range("H25").Value = 143.543553
$ws.Range("I25").Value = 0.1574303305363417
$ws.Range("J25").Value = 0.1574303305363417
$ws.Range("M25").Value = 5.865492
$ws.Range("N25").Value = 17.596476
$ws.Range("O25").Value = 0.2208713818117337
$ws.Range("P25").Value = 0.2208713818117337
$ws.Range("Q25").Value = 280.651187257692
$ws.Range("R25").Value = 2525.860685319228
$ws.Range("S25").Value = 0.03477185464463977
$ws.Range("T25").Value = 0.03477185464463978
$ws.Range("G26").Value = 47.847851
$ws.Range("H26").Value = 143.543553
$ws.Range("I26").Value = 0.1574303305363417
$ws.Range("J26").Value = 0.1574303305363417
$ws.Range("M26").Value = 2.384473
$ws.Range("N26").Value = 7.153419
$ws.Range("O26").Value = 0.08978988402043174
$ws.Range("P26").Value = 0.08978988402043174
$ws.Range("Q26").Value = 114.091908817523
$ws.Range("R26").Value = 1026.827179357707
$ws.Range("S26").Value = 0.01413565112015636
$ws.Range("T26").Value = 0.01413565112015636
